$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L2").Value = 0.9624999761581421
$ws.Range("M2").Value = 0.9833333492279053
$ws.Range("M3").Value = 0.8999999761581421
$ws.Range("L4").Value = 0.9895833134651184
$ws.Range("M4").Value = 0.8833333253860474
$ws.Range("L5").Value = 0.987500011920929
$ws.Range("M5").Value = 0.9666666388511658
$ws.Range("L6").Value = 0.9854166507720947
$ws.Range("L7").Value = 0.9791666865348816
$ws.Range("M7").Value = 1
$ws.Range("L8").Value = 0.9708333611488342
$ws.Range("M8").Value = 0.8666666746139526
$ws.Range("M9").Value = 0.8833333253860474
$ws.Range("L10").Value = 0.9895833134651184
$ws.Range("M10").Value = 0.9166666865348816
$ws.Range("L11").Value = 0.9770833253860474
$ws.Range("M11").Value = 0.8999999761581421
$ws.Range("L12").Value = 0.9854166507720947
$ws.Range("M12").Value = 0.8833333253860474
$ws.Range("L13").Value = 0.9854166507720947
$ws.Range("M13").Value = 0.9333333373069763
$ws.Range("L14").Value = 0.9895833134651184
$ws.Range("M14").Value = 0.9833333492279053
$ws.Range("L15").Value = 0.9895833134651184
$ws.Range("M15").Value = 0.949999988079071
$ws.Range("L16").Value = 0.9895833134651184
$ws.Range("M16").Value = 0.8166666626930237
$ws.Range("L18").Value = 0.987500011920929
$ws.Range("M18").Value = 1
$ws.Range("L19").Value = 0.9916666746139526
$ws.Range("M19").Value = 0.8500000238418579
$ws.Range("L20").Value = 0.9791666865348816
$ws.Range("M20").Value = 0.9333333373069763
$ws.Range("L21").Value = 0.987500011920929
$ws.Range("L22").Value = 0.981249988079071
$ws.Range("M22").Value = 0.9333333373069763
$ws.Range("L23").Value = 0.9708333611488342
$ws.Range("M23").Value = 0.949999988079071
$ws.Range("L24").Value = 0.9791666865348816
$ws.Range("L25").Value = 0.981249988079071
$ws.Range("M25").Value = 0.9833333492279053
$ws.Range("L26").Value = 0.9729166626930237
$ws.Range("M26").Value = 0.949999988079071
$ws.Range("L27").Value = 0.987500011920929
$ws.Range("M27").Value = 0.9333333373069763
$ws.Range("L28").Value = 0.987500011920929
$ws.Range("M28").Value = 0.6333333253860474
$ws.Range("L29").Value = 0.9937499761581421
$ws.Range("M29").Value = 0.9833333492279053
$ws.Range("L30").Value = 0.987500011920929
$ws.Range("M30").Value = 0.9833333492279053
$ws.Range("L31").Value = 0.987500011920929
$ws.Range("M31").Value = 0.8166666626930237
$ws.Range("L32").Value = 0.9770833253860474
$ws.Range("M32").Value = 0.8999999761581421
$ws.Range("L33").Value = 0.9729166626930237
$ws.Range("M33").Value = 0.8999999761581421
$ws.Range("L34").Value = 0.987500011920929
$ws.Range("M34").Value = 0.9833333492279053
$ws.Range("L35").Value = 0.9916666746139526
$ws.Range("M35").Value = 0.9333333373069763
$ws.Range("L36").Value = 0.9833333492279053
$ws.Range("L37").Value = 0.987500011920929
$ws.Range("M37").Value = 0.9666666388511658
$ws.Range("L39").Value = 0.9833333492279053
$ws.Range("L40").Value = 0.9916666746139526
$ws.Range("M40").Value = 0.8833333253860474
$ws.Range("L41").Value = 0.9854166507720947
$ws.Range("L42").Value = 0.9895833134651184
$ws.Range("L43").Value = 1
$ws.Range("M43").Value = 0.8999999761581421
$ws.Range("L44").Value = 0.9604166746139526
$ws.Range("M44").Value = 0.8333333134651184
$ws.Range("M45").Value = 0.8833333253860474
$ws.Range("L46").Value = 0.9833333492279053
$ws.Range("M46").Value = 0.6833333373069763
$ws.Range("L47").Value = 0.96875
$ws.Range("M47").Value = 0.949999988079071
$ws.Range("L48").Value = 0.9854166507720947
$ws.Range("M48").Value = 0.9833333492279053
$ws.Range("L49").Value = 0.9979166388511658
$ws.Range("L50").Value = 0.9770833253860474
$ws.Range("M50").Value = 1
$ws.Range("L51").Value = 0.9791666865348816
$ws.Range("M51").Value = 0.9833333492279053
$ws.Range("L52").Value = 0.981249988079071
$ws.Range("M52").Value = 0.8999999761581421
$ws.Range("L53").Value = 0.981249988079071
$ws.Range("M53").Value = 0.9833333492279053
$ws.Range("L54").Value = 0.9958333373069763
$ws.Range("M54").Value = 0.9666666388511658
$ws.Range("M55").Value = 0.949999988079071
$ws.Range("L56").Value = 0.9729166626930237
$ws.Range("M56").Value = 0.949999988079071
$ws.Range("L57").Value = 0.9750000238418579
$ws.Range("M57").Value = 1
$ws.Range("L58").Value = 0.9833333492279053
$ws.Range("M58").Value = 0.949999988079071
$ws.Range("L59").Value = 0.987500011920929
$ws.Range("M59").Value = 1
$ws.Range("L60").Value = 0.987500011920929
$ws.Range("M60").Value = 0.949999988079071
$ws.Range("M61").Value = 0.699999988079071
$ws.Range("L62").Value = 0.9791666865348816
$ws.Range("L63").Value = 0.981249988079071
$ws.Range("M63").Value = 0.8666666746139526
$ws.Range("L64").Value = 0.9916666746139526
$ws.Range("M64").Value = 0.5666666626930237
$ws.Range("L65").Value = 0.987500011920929
$ws.Range("M65").Value = 0.8999999761581421
$ws.Range("L66").Value = 0.9937499761581421
$ws.Range("M66").Value = 0.9333333373069763
$ws.Range("L67").Value = 0.9937499761581421
$ws.Range("M67").Value = 0.8999999761581421
$ws.Range("L68").Value = 0.9624999761581421
$ws.Range("M68").Value = 0.9333333373069763
$ws.Range("L69").Value = 0.981249988079071
$ws.Range("M69").Value = 0.9833333492279053
$ws.Range("M70").Value = 0.8999999761581421
$ws.Range("L71").Value = 0.9854166507720947
$ws.Range("M71").Value = 0.9166666865348816
$ws.Range("L72").Value = 0.9937499761581421
$ws.Range("L73").Value = 0.9750000238418579
$ws.Range("M73").Value = 0.8833333253860474
$ws.Range("L74").Value = 0.96875
$ws.Range("M74").Value = 0.949999988079071
$ws.Range("L75").Value = 0.9937499761581421
$ws.Range("M75").Value = 0.8999999761581421
$ws.Range("L76").Value = 0.9937499761581421
$ws.Range("M76").Value = 0.8333333134651184
$ws.Range("L77").Value = 0.96875
$ws.Range("M77").Value = 0.949999988079071
$ws.Range("L78").Value = 0.9937499761581421
$ws.Range("M78").Value = 0.949999988079071
$ws.Range("M79").Value = 0.9166666865348816
$ws.Range("L80").Value = 0.9458333253860474
$ws.Range("L81").Value = 0.96875
$ws.Range("M81").Value = 0.9666666388511658
$ws.Range("L82").Value = 0.981249988079071
$ws.Range("M82").Value = 0.9666666388511658
$ws.Range("L83").Value = 0.9729166626930237
$ws.Range("M83").Value = 0.9833333492279053
$ws.Range("L84").Value = 0.9895833134651184
$ws.Range("M84").Value = 0.9833333492279053
$ws.Range("L85").Value = 0.9958333373069763
$ws.Range("M85").Value = 0.9666666388511658
$ws.Range("L86").Value = 0.9645833373069763
$ws.Range("M86").Value = 0.9666666388511658
$ws.Range("L87").Value = 0.9854166507720947
$ws.Range("M87").Value = 0.9833333492279053
$ws.Range("L88").Value = 0.9895833134651184
$ws.Range("L89").Value = 0.9770833253860474
$ws.Range("M90").Value = 0.9666666388511658
$ws.Range("L91").Value = 0.9979166388511658
$ws.Range("M91").Value = 0.9833333492279053
$ws.Range("L92").Value = 0.9645833373069763
$ws.Range("M92").Value = 0.949999988079071
$ws.Range("L93").Value = 0.9791666865348816
$ws.Range("M93").Value = 0.9833333492279053
$ws.Range("L94").Value = 0.9791666865348816
$ws.Range("M94").Value = 0.9666666388511658
$ws.Range("L95").Value = 0.9895833134651184
$ws.Range("L96").Value = 0.981249988079071
$ws.Range("M96").Value = 1
